$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-31 from
# 2023-09-05 (45174) to 2023-09-06 (45175)
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 3).Value = 45175
}
